# Update scoring data for TongHopDiem1 sheet (rows 4-8, columns A:AQ)
# Replaces placeholder rows with real student names/IDs and committee scores
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowData = @{
    4 = @("La Thị Hồng Nhung", "B20DCDT158", "D20XLTH", 8, 8, 8, 8, 9, 9, 9, 9, 5, 5, 5, 5, 6, 6, 6, 6, 1, 1, 1, 1, 7, 7, 8, 8, 8, 9, 9, 9, 9.5, 7, 8, 7, 7.5, 7, 8, 6.75, 6.428571428571429, 6.285714285714286, 6.571428571428571, 6.839285714285714)
    5 = @("Nguyễn Trung Tuấn", "B20DCDT190", "D20XLTH", 8, 8, 8, 8, 9, 9, 9, 9, 5, 5, 5, 5, 6, 6, 6, 6, 1, 1, 1, 1, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 3.625, 4.142857142857143, 4.142857142857143, 4.142857142857143, 2.675595238095239)
    6 = @("Phạm Thế Anh", "B20DCDT017", "D20DTMT1", 3, 3, 3, 3, 4, 4, 4, 4, 5, 5, 5, 5, 6, 6, 6, 6, 7, 7, 7, 7, 1, 1, 0, 0, 0, 0, 0, 0, 0, 2, 2, 2, 2, 1, 0.6666666666666666, 3.375, 3.857142857142857, 3.714285714285714, 3.857142857142857, 2.745039682539682)
    7 = @("Nguyễn Tiến Duy", "B20DCDT037", "D20DTMT1", 3, 3, 3, 3, 4, 4, 4, 4, 5, 5, 5, 5, 6, 6, 6, 6, 7, 7, 7, 7, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 3.125, 3.571428571428572, 3.571428571428572, 3.571428571428572, 2.306547619047619)
    8 = @("Lê Sỹ Sang", "B20DCDT175", "D20DTMT2", 3, 3, 3, 3, 4, 4, 4, 4, 5, 5, 5, 5, 6, 6, 6, 6, 7, 7, 7, 7, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 3.125, 3.571428571428572, 3.571428571428572, 3.571428571428572, 2.306547619047619)
}

foreach ($r in $rowData.Keys) {
    $vals = $rowData[$r]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $ws.Cells.Item($r, $i + 1).Value = $vals[$i]
    }
}
